# Weekly update: shift existing Pera price rows down by two (new week's
# records are inserted at the top of the data block, rows 100-101), and
# append the two rows that fall off the bottom of the original range as
# new rows 133-134.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 100:132 down to 102:134, preserving values and formatting
# (Excel resolves this from a snapshot of the source range, so the
# overlap between source and destination is handled correctly).
$src = $ws.Range("A100:T132")
$dst = $ws.Range("A102:T134")
$src.Copy($dst)

# New week's entries for row 100 (Especial) and row 101 (Primera)
$ws.Range("D100").Value = 44466
$ws.Range("L100").Value = "Especial"
$ws.Range("M100").Value = 30
$ws.Range("N100").Value = 11000
$ws.Range("O100").Value = 11000
$ws.Range("P100").Value = 11000
$ws.Range("S100").Value = 688

$ws.Range("D101").Value = 44466
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 60
$ws.Range("N101").Value = 9000
$ws.Range("O101").Value = 10000
$ws.Range("P101").Value = 9500
$ws.Range("S101").Value = 594
